$wb = $excel.ActiveWorkbook
$players = $wb.Worksheets.Item("Players")
$owners = $wb.Worksheets.Item("OwnerTotals")

# --- Players sheet (sheet1) cell updates ---
$players.Range("G7").Value = "3:32 - 2nd Half"
$players.Range("O7").Value = 34
$players.Range("G8").Value = "1:40 - 2nd Half"
$players.Range("H8").Value = 11
$players.Range("I8").Value = 16
$players.Range("O8").Value = 32
$players.Range("G9").Value = "7:57 - 2nd Half"
$players.Range("O9").Value = 24
$players.Range("G10").Value = "7:57 - 2nd Half"
$players.Range("H10").Value = -1
$players.Range("O10").Value = 20
$players.Range("G11").Value = "1:40 - 2nd Half"
$players.Range("G16").Value = "1:40 - 2nd Half"
$players.Range("O16").Value = 35
$players.Range("G17").Value = "7:57 - 2nd Half"
$players.Range("O17").Value = 20
$players.Range("G18").Value = "1:40 - 2nd Half"
$players.Range("G19").Value = "3:32 - 2nd Half"
$players.Range("D27").Value = "Keyshawn Hall"
$players.Range("E27").Value = "AUB"
$players.Range("F27").Value = "TEX@AUB"
$players.Range("G27").Value = "7:57 - 2nd Half"
$players.Range("H27").Value = 14
$players.Range("I27").Value = 20
$players.Range("J27").Value = 1
$players.Range("K27").Value = 2
$players.Range("L27").Value = 0
$players.Range("N27").Value = 2
$players.Range("O27").Value = 30
$players.Range("D28").Value = "Dedan Thomas Jr."
$players.Range("E28").Value = "LSU"
$players.Range("F28").Value = "MSST@LSU"
$players.Range("G28").Value = "3:32 - 2nd Half"
$players.Range("H28").Value = 13
$players.Range("I28").Value = 14
$players.Range("J28").Value = 2
$players.Range("K28").Value = 4
$players.Range("L28").Value = 1
$players.Range("N28").Value = 0
$players.Range("O28").Value = 28
$players.Range("G29").Value = "3:32 - 2nd Half"
$players.Range("H29").Value = 19
$players.Range("I29").Value = 13
$players.Range("O29").Value = 22
$players.Range("G36").Value = "3:32 - 2nd Half"
$players.Range("G37").Value = "1:40 - 2nd Half"
$players.Range("G42").Value = "1:40 - 2nd Half"
$players.Range("H42").Value = 9
$players.Range("J42").Value = 3
$players.Range("G43").Value = "3:32 - 2nd Half"
$players.Range("O43").Value = 15
$players.Range("G44").Value = "7:57 - 2nd Half"
$players.Range("H44").Value = 19
$players.Range("J44").Value = 1
$players.Range("O44").Value = 26
$players.Range("G45").Value = "1:40 - 2nd Half"
$players.Range("O45").Value = 26
$players.Range("G46").Value = "7:57 - 2nd Half"
$players.Range("O46").Value = 28
$players.Range("G50").Value = "7:57 - 2nd Half"
$players.Range("H50").Value = 25
$players.Range("I50").Value = 24
$players.Range("O50").Value = 31
$players.Range("G51").Value = "1:40 - 2nd Half"
$players.Range("G52").Value = "1:40 - 2nd Half"
$players.Range("G53").Value = "3:32 - 2nd Half"
$players.Range("H53").Value = 9
$players.Range("J53").Value = 2
$players.Range("O53").Value = 31
$players.Range("G59").Value = "7:57 - 2nd Half"
$players.Range("G60").Value = "3:32 - 2nd Half"
$players.Range("H60").Value = 11
$players.Range("O60").Value = 27
$players.Range("G61").Value = "3:32 - 2nd Half"
$players.Range("O61").Value = 26
$players.Range("G62").Value = "7:57 - 2nd Half"
$players.Range("O62").Value = 30
$players.Range("G90").Value = "1:40 - 2nd Half"
$players.Range("O90").Value = 26
$players.Range("G91").Value = "3:32 - 2nd Half"
$players.Range("O91").Value = 26
$players.Range("G92").Value = "1:40 - 2nd Half"
$players.Range("G93").Value = "7:57 - 2nd Half"
$players.Range("O93").Value = 31
$players.Range("G94").Value = "1:40 - 2nd Half"
$players.Range("H94").Value = 12
$players.Range("J94").Value = 1
$players.Range("G95").Value = "3:32 - 2nd Half"
$players.Range("H95").Value = 11
$players.Range("J95").Value = 3
$players.Range("M95").Value = 1
$players.Range("O95").Value = 26
$players.Range("G96").Value = "3:32 - 2nd Half"
$players.Range("G97").Value = "3:32 - 2nd Half"
$players.Range("G98").Value = "1:40 - 2nd Half"
$players.Range("G99").Value = "1:40 - 2nd Half"
$players.Range("D100").Value = "Sebastian Williams-Adams"
$players.Range("E100").Value = "AUB"
$players.Range("F100").Value = "TEX@AUB"
$players.Range("G100").Value = "7:57 - 2nd Half"
$players.Range("I100").Value = 7
$players.Range("K100").Value = 1
$players.Range("L100").Value = 2
$players.Range("M100").Value = 2
$players.Range("N100").Value = 0
$players.Range("O100").Value = 26
$players.Range("D101").Value = "Rashad King"
$players.Range("E101").Value = "LSU"
$players.Range("F101").Value = "MSST@LSU"
$players.Range("G101").Value = "3:32 - 2nd Half"
$players.Range("H101").Value = 7
$players.Range("I101").Value = 11
$players.Range("K101").Value = 0
$players.Range("L101").Value = 1
$players.Range("M101").Value = 0
$players.Range("N101").Value = 2
$players.Range("O101").Value = 18
$players.Range("G102").Value = "1:40 - 2nd Half"
$players.Range("G103").Value = "7:57 - 2nd Half"
$players.Range("G104").Value = "3:32 - 2nd Half"
$players.Range("O104").Value = 17
$players.Range("G105").Value = "3:32 - 2nd Half"
$players.Range("G106").Value = "1:40 - 2nd Half"
$players.Range("G107").Value = "1:40 - 2nd Half"
$players.Range("G108").Value = "7:57 - 2nd Half"
$players.Range("G109").Value = "1:40 - 2nd Half"
$players.Range("G110").Value = "7:57 - 2nd Half"
$players.Range("G111").Value = "1:40 - 2nd Half"
$players.Range("G112").Value = "3:32 - 2nd Half"
$players.Range("O112").Value = 22
$players.Range("G113").Value = "7:57 - 2nd Half"
$players.Range("G114").Value = "1:40 - 2nd Half"
$players.Range("G115").Value = "7:57 - 2nd Half"
$players.Range("G116").Value = "3:32 - 2nd Half"
$players.Range("H51").Value = 10
$players.Range("I51").Value = 8

# --- Column G width update (18 -> 17) ---
$players.Columns.Item(7).ColumnWidth = 16.1

# --- OwnerTotals sheet (sheet2) cell updates ---
$owners.Range("B2").Value = 61
$owners.Range("B3").Value = 52
$owners.Range("A4").Value = "Ron"
$owners.Range("C4").Value = 3
$owners.Range("A6").Value = "Mark"
$owners.Range("B6").Value = 45
$owners.Range("C6").Value = 4
